$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.643.31'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").Value = '3.391.03'
$ws.Range("E3").Value = '  +1.70%  '

$ws.Range("D4").Value = "'" + '0.999'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = "'" + '183.94'
$ws.Range("E5").Value = '  +1.21%  '

$ws.Range("D6").Value = "'" + '540.44'
$ws.Range("E6").Value = '  +1.34%  '

$ws.Range("D7").Value = "'" + '0.604'
$ws.Range("E7").Value = '  -0.60%  '

$ws.Range("D8").Value = '3.384.82'
$ws.Range("E8").Value = '  +1.69%  '

$ws.Range("D9").Value = "'" + '1.00'
$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = "'" + '0.627'
$ws.Range("E10").Value = '  +1.82%  '

$ws.Range("D11").Value = "'" + '55.95'
$ws.Range("E11").Value = '  -6.47%  '

$ws.Range("D12").Value = "'" + '0.141'
$ws.Range("E12").Value = '  +4.48%  '

$ws.Range("D13").Value = "'" + '0.0000268'
$ws.Range("E13").Value = '  +2.37%  '

$ws.Range("D14").Value = "'" + '9.27'
$ws.Range("E14").Value = '  +0.81%  '

$ws.Range("D15").Value = '3.912.51'
$ws.Range("E15").Value = '  +1.94%  '

$ws.Range("E16").Value = '  +1.73%  '

$ws.Range("D17").Value = '3.375.85'
$ws.Range("E17").Value = '  +1.80%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = "'" + '18.00'
$ws.Range("E18").Value = '  +1.73%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '65.773.61'
$ws.Range("E19").Value = '  +1.65%  '

$ws.Range("D20").Value = "'" + '11.46'
$ws.Range("E20").Value = '  +1.41%  '

$ws.Range("D21").Value = "'" + '0.989'
$ws.Range("E21").Value = '  +2.11%  '

$ws.Range("D22").Value = "'" + '391.18'
$ws.Range("E22").Value = '  +3.38%  '

$ws.Range("D23").Value = "'" + '12.22'
$ws.Range("E23").Value = '  +7.70%  '

$ws.Range("D24").Value = "'" + '4.22'
$ws.Range("E24").Value = '  +6.74%  '

$ws.Range("D25").Value = "'" + '83.17'
$ws.Range("E25").Value = '  +2.38%  '

$ws.Range("D26").Value = "'" + '3.81'
$ws.Range("E26").Value = '  -0.91%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = "'" + '2.84'
$ws.Range("E27").Value = '  +4.85%  '

$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").Value = "'" + '6.12'
$ws.Range("E28").Value = '  +0.35%  '

$ws.Range("D29").Value = "'" + '11.58'
$ws.Range("E29").Value = '  -0.54%  '

$ws.Range("D30").Value = "'" + '8.51'
$ws.Range("E30").Value = '  +0.31%  '

$ws.Range("D31").Value = "'" + '29.62'
$ws.Range("E31").Value = '  +1.30%  '

$ws.Range("D32").Value = "'" + '666.38'
$ws.Range("E32").Value = '  +0.73%  '

$ws.Range("D33").Value = "'" + '6.87'
$ws.Range("E33").Value = '  +1.65%  '

$ws.Range("D34").Value = "'" + '11.49'
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("E35").Value = '  +1.41%  '

$ws.Range("D36").Value = "'" + '58.00'
$ws.Range("E36").Value = '  -3.14%  '

$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").Value = "'" + '38.10'
$ws.Range("E37").Value = '  +2.47%  '

$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D38").Value = "'" + '0.400'
$ws.Range("E38").Value = '  +0.47%  '

$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = "'" + '1.00'
$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("E40").Value = '  +10.77%  '

$ws.Range("D41").Value = "'" + '2.83'
$ws.Range("E41").Value = '  +11.16%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'" + '3.29'
$ws.Range("E42").Value = '  +15.65%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = "'" + '0.997'
$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("D44").Value = "'" + '0.130'
$ws.Range("E44").Value = '  +1.58%  '

$ws.Range("D45").Value = '3.033.82'
$ws.Range("E45").Value = '  +3.45%  '

$ws.Range("D46").Value = "'" + '2.77'
$ws.Range("E46").Value = '  +1.19%  '

$ws.Range("D47").Value = "'" + '0.0413'
$ws.Range("E47").Value = '  +2.29%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = "'" + '2.70'
$ws.Range("E48").Value = '  +1.30%  '

$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = "'" + '3.15'
$ws.Range("E49").Value = '  +0.78%  '

$ws.Range("D50").Value = "'" + '0.128'
$ws.Range("E50").Value = '  +0.19%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = "'" + '8.55'
$ws.Range("E51").Value = '  +6.68%  '
